# Apply "Penalty Reward System" forecast shift edit.
$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Sheet: Forecast Comparison ---
# Shift Week_Start_Date (col B) forward by one week (7 days) for rows 2-17,
# and zero out the MyForecast (col D) values for those rows.
$dates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $cell = $wsForecast.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$i]
    $wsForecast.Cells.Item($row, 4).Value = 0
}

# --- Sheet: Summary ---
$summaryUpdates = @{
    "B2"  = "2022-12-25 to 2025-01-05"
    "B9"  = "3"
    "B10" = "2"
    "B11" = "1"
    "B12" = "0"
    "B13" = "2025-01-12"
    "B14" = "0"
}

foreach ($addr in $summaryUpdates.Keys) {
    $cell = $wsSummary.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $summaryUpdates[$addr]
}
